# Results.xlsx - "Finished discussion, conclusion and future work"
#
# The author picked a different scenario in the MasterSheet dropdown (J1),
# which re-pointed the SWITCH() lookup from "400 vs 400" to "800 vs 800"
# and recalculated the summary table (B2:F18). They also navigated the
# workbook a bit while doing so, which moved the selection / active-tab
# state on a few sheets (MasterSheet becomes the active tab, the
# 12800vs12800 tab is no longer active, and the Charts sheet scrolled to a
# different selection).

$wb = $excel.ActiveWorkbook

# --- Charts sheet: update the lingering selection (not the active tab) ---
$wsCharts = $wb.Worksheets.Item("Charts")
$wsCharts.Activate() | Out-Null
$wsCharts.Range("M171").Select() | Out-Null

# --- 12800vs12800 sheet: update its lingering selection; it stops being
#     the active/tabSelected sheet once we move on ---
$ws12800 = $wb.Worksheets.Item("12800vs12800")
$ws12800.Activate() | Out-Null
$ws12800.Range("I8").Select() | Out-Null

# --- MasterSheet: change the scenario dropdown, which recalculates the
#     SWITCH-based summary formulas, and becomes the final active sheet ---
$wsMaster = $wb.Worksheets.Item("MasterSheet")
$wsMaster.Activate() | Out-Null
$wsMaster.Range("J1").Value = "800 vs 800"
$wsMaster.Range("I6").Select() | Out-Null
